{"js": "// The document has one bookmark named \"smarthosting\" whose w:id needs to\n// change from \"1\" to \"0\" (Office.js does not expose raw bookmark ids, but\n// deleting + re-inserting a lone bookmark causes the engine to re-mint it\n// starting at id 0, which is exactly the target state).\n//\n// It also has a paragraph (\"100,000 Smart: Arguably the hardest part. Smart\n// can be obtained from exchanges such as CryptoBridge, HitBTC. For the full\n// list ...\") that loses the two HYPERLINK field-coded mentions of\n// \"CryptoBridge\" and \"HitBTC\" (including their begin/instrText/separate/end\n// runs and the spell-check proofErr markers) as well as the now-dangling\n// \" such as\" / \", \" connector text, leaving \"...exchanges. For the full\n// list...\".\n\nconst body = context.document.body;\n\n// --- 1. Re-mint the \"smarthosting\" bookmark so its id becomes 0 ---------\ncontext.document.deleteBookmark(\"smarthosting\");\nawait context.sync();\n\nconst headingParas = body.paragraphs;\nheadingParas.load(\"items/text\");\nawait context.sync();\n\nconst headingPara = headingParas.items.find(\n  (p) => p.text === \"WHAT DO I NEED TO HOST A SMARTNODE?\"\n);\nif (!headingPara) {\n  throw new Error('Could not find the \"WHAT DO I NEED TO HOST A SMARTNODE?\" paragraph');\n}\nconst headingStart = headingPara.getRange(\"Start\");\nheadingStart.insertBookmark(\"smarthosting\");\nawait context.sync();\n\n// --- 2. Strip the CryptoBridge / HitBTC exchange mentions ----------------\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst targetPara = paras.items.find(\n  (p) => p.text.indexOf(\"Arguably the hardest part\") !== -1\n);\nif (!targetPara) {\n  throw new Error(\"Could not find the SmartNode hosting requirements paragraph\");\n}\n\n// Remove the two HYPERLINK field codes (CryptoBridge, HitBTC) -- this also\n// removes their begin/instrText/separate/result/end runs and proofErr marks.\nconst fields = targetPara.fields;\nfields.load(\"items/code\");\nawait context.sync();\n\nconst fieldsToDelete = fields.items.filter((f) => f.code.indexOf(\"HYPERLINK\") !== -1 &&\n  (f.code.indexOf(\"crypto-bridge.org\") !== -1 || f.code.indexOf(\"hitbtc.com\") !== -1));\nfieldsToDelete.forEach((f) => f.delete());\nawait context.sync();\n\n// Remove the now-dangling \" such as<nbsp>,<nbsp>\" connector text that used\n// to glue \"exchanges\" to the two links, leaving just \"...exchanges.\".\nconst suchAsHit = targetPara.search(\" such as\", { matchCase: false, matchWholeWord: false });\nsuchAsHit.load(\"items\");\nconst commaHit = targetPara.search(\",\\u00a0\", { matchCase: false, matchWholeWord: false });\ncommaHit.load(\"items\");\nawait context.sync();\n\nif (suchAsHit.items.length > 0 && commaHit.items.length > 0) {\n  const toRemove = suchAsHit.items[0].expandTo(commaHit.items[0]);\n  toRemove.delete();\n  await context.sync();\n}\n", "ps1": "# The document has one bookmark named \"smarthosting\" whose w:id needs to\n# change from \"1\" to \"0\". Word COM does not expose the raw bookmark id\n# directly, so we recreate the bookmark (delete + re-Add) at the exact same\n# (collapsed, zero-length) location -- the start of the\n# \"WHAT DO I NEED TO HOST A SMARTNODE?\" heading paragraph -- which is where\n# the engine mints the next available id (0, since it's the only bookmark).\n#\n# It also has a paragraph (\"100,000 Smart: Arguably the hardest part. Smart\n# can be obtained from exchanges such as CryptoBridge, HitBTC. For the full\n# list ...\") that loses the two HYPERLINK field-coded mentions of\n# \"CryptoBridge\" and \"HitBTC\" (their begin/instrText/separate/result/end\n# runs and the spell-check proofErr markers) as well as the now-dangling\n# \" such as\" / \", \" connector text, leaving \"...exchanges. For the full\n# list...\".\n\n$d = $word.ActiveDocument\n\n# --- 1. Re-mint the \"smarthosting\" bookmark so its id becomes 0 -----------\n$headingPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"WHAT DO I NEED TO HOST A SMARTNODE?\") {\n        $headingPara = $p\n        break\n    }\n}\nif ($headingPara -eq $null) {\n    throw \"Could not find the 'WHAT DO I NEED TO HOST A SMARTNODE?' paragraph\"\n}\n$bmRange = $headingPara.Range.Duplicate\n$bmRange.Collapse(1)  # wdCollapseStart\n\nif ($d.Bookmarks.Exists(\"smarthosting\")) {\n    $d.Bookmarks.Item(\"smarthosting\").Delete()\n}\n$d.Bookmarks.Add(\"smarthosting\", $bmRange)\n\n# --- 2. Strip the CryptoBridge / HitBTC exchange mentions ------------------\n# Remove the two HYPERLINK field codes -- this also removes their\n# begin/instrText/separate/result/end runs and proofErr marks.\nfor ($i = $d.Fields.Count; $i -ge 1; $i--) {\n    $f = $d.Fields.Item($i)\n    $code = $f.Code.Text\n    if ($code -like \"*crypto-bridge.org*\" -or $code -like \"*hitbtc.com*\") {\n        $f.Delete()\n    }\n}\n\n# Find the target paragraph (still paragraph-stable after the field deletes).\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Arguably the hardest part*\") {\n        $targetPara = $p\n        break\n    }\n}\nif ($targetPara -eq $null) {\n    throw \"Could not find the SmartNode hosting requirements paragraph\"\n}\n\n# Remove the now-dangling \" such as<nbsp>,<nbsp>\" connector text that used\n# to glue \"exchanges\" to the two links, leaving just \"...exchanges.\".\n$nbsp = [char]0x00A0\n$searchRange = $targetPara.Range.Duplicate\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \" such as\" + $nbsp + \",\" + $nbsp\n$find.Replacement.Text = \"\"\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop\n\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n"}
